$wb = $excel.ActiveWorkbook

# --- Insert the new "daily" sheet, positioned right before "mouth" so it
#     lands between "pair" and "mouth" (sheetId 7 / 3rd tab overall). ---
$mouthSheet = $wb.Worksheets("mouth")
$daily = $wb.Worksheets.Add($mouthSheet)
$daily.Name = "daily"

# Header row matching the other "blank" vocab sheets (hand / sense):
# French | English | French | English
$daily.Range("A1").Value = "French"
$daily.Range("B1").Value = "English"
$daily.Range("C1").Value = "French"
$daily.Range("D1").Value = "English"
$daily.Range("A1:D1").HorizontalAlignment = -4108

# View: new sheet becomes the active tab, zoomed to 173%, A1:D1 selected.
$daily.Activate()
$excel.ActiveWindow.Zoom = 173
$daily.Range("A1:D1").Select()

# --- Clean up the duplicated "center, no fill" style on the "mouth" sheet
#     so those cells share the same style as the rest (instead of the
#     redundant twin style). ---
$mouth = $wb.Worksheets("mouth")
$mouth.Range("D3").HorizontalAlignment = -4108
$mouth.Range("A5:B5").HorizontalAlignment = -4108
$mouth.Range("A9:B9").HorizontalAlignment = -4108
$mouth.Range("A11:B11").HorizontalAlignment = -4108
$mouth.Range("A13:B13").HorizontalAlignment = -4108

# --- Same cleanup on the "leg" sheet. ---
$leg = $wb.Worksheets("leg")
$leg.Range("C2:D2").HorizontalAlignment = -4108
$leg.Range("A6").HorizontalAlignment = -4108
